# Horarios Linea 141 - scrape update 16:28:03 -> 16:37:06
# Applies the new-scrape diff across the three worksheets:
#   LP1912       (sheet1) - reorders + 4 new rows (2 inserted mid-table, 2 appended)
#   LP1912-215   (sheet2) - timestamp refresh only
#   6203-6173    (sheet3) - timestamp refresh + 1 appended row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header block
$ws1.Cells.Item(2,1).Value = "Última actualización: 16:37:06"
$ws1.Cells.Item(3,1).Value = "Total filas: 265"

# Rows 86 / 87 swap (Hora_Scrap, Linea, Minutos) while Hora_Llegada/Parada stay put
$ws1.Cells.Item(86,1).Value = "09:23:23"
$ws1.Cells.Item(86,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(86,4).Value = 96

$ws1.Cells.Item(87,1).Value = "10:50:41"
$ws1.Cells.Item(87,3).Value = "10_OLMOS"
$ws1.Cells.Item(87,4).Value = 9

# Rows 120 / 121 swap (Linea only)
$ws1.Cells.Item(120,3).Value = "15_ABASTO"
$ws1.Cells.Item(121,3).Value = "16_P MOR-SANTA ANA"

# Rows 137 / 138 swap (Hora_Scrap, Linea, Minutos)
$ws1.Cells.Item(137,1).Value = "11:52:01"
$ws1.Cells.Item(137,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(137,4).Value = 45

$ws1.Cells.Item(138,1).Value = "11:47:17"
$ws1.Cells.Item(138,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(138,4).Value = 50

# Rows 169 / 170 swap (Hora_Scrap, Linea, Minutos)
$ws1.Cells.Item(169,1).Value = "12:33:21"
$ws1.Cells.Item(169,3).Value = "10_OLMOS"
$ws1.Cells.Item(169,4).Value = 89

$ws1.Cells.Item(170,1).Value = "13:14:29"
$ws1.Cells.Item(170,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(170,4).Value = 48

# Two freshly-scraped rows land right after row 256 ("17:49" / "17:53" arrivals
# split), pushing the previous rows 257-266 down to 259-268 unchanged.
$ws1.Rows.Item(257).Resize(2).Insert()

$ws1.Cells.Item(257,1).Value = "16:37:06"
$ws1.Cells.Item(257,2).Value = "17:50"
$ws1.Cells.Item(257,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(257,4).Value = 73
$ws1.Cells.Item(257,5).Value = "LP1912"

$ws1.Cells.Item(258,1).Value = "16:37:06"
$ws1.Cells.Item(258,2).Value = "17:52"
$ws1.Cells.Item(258,3).Value = "10_OLMOS"
$ws1.Cells.Item(258,4).Value = 75
$ws1.Cells.Item(258,5).Value = "LP1912"

# Two brand-new rows appended at the bottom of the table
$ws1.Cells.Item(269,1).Value = "16:37:06"
$ws1.Cells.Item(269,2).Value = "18:30"
$ws1.Cells.Item(269,3).Value = "14_ABASTO"
$ws1.Cells.Item(269,4).Value = 113
$ws1.Cells.Item(269,5).Value = "LP1912"

$ws1.Cells.Item(270,1).Value = "16:37:06"
$ws1.Cells.Item(270,2).Value = "18:36"
$ws1.Cells.Item(270,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(270,4).Value = 119
$ws1.Cells.Item(270,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 16:37:06"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 16:37:06"
$ws3.Cells.Item(3,1).Value = "Total filas: 36"

$ws3.Cells.Item(41,1).Value = "16:37:06"
$ws3.Cells.Item(41,2).Value = "18:36"
$ws3.Cells.Item(41,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(41,4).Value = 119
$ws3.Cells.Item(41,5).Value = "L6203"
